# Atualização automática de CACHOEIRINHA.xlsx
$wb = $excel.ActiveWorkbook

# Delete the "Desarquivamentos Pendentes" sheet entirely (it has no other
# sheet depending on its data).
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()
$excel.DisplayAlerts = $true

# Rename the remaining sheets to their new (uppercased) titles.
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"
